# Generate Report for Handback
# Adds a new handback entry (d452f1dd-74d8-4c8f-972c-7be67665c439) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets, expanding each sheet's
# table by one row and wiring up the matching hyperlinks.

$wb = $excel.ActiveWorkbook

$fileBase   = "d452f1dd-74d8-4c8f-972c-7be67665c439"
$mdName     = "$fileBase.md"
$mdPath     = "e2e\$fileBase.md"
$blobHash   = "b07018d703e746de6f5096f3692d55169f1d988e"
$zhXlf      = "$fileBase.$blobHash.zh-cn.xlf"
$deXlf      = "$fileBase.$blobHash.de-de.xlf"

$zhHandoffDt  = "2016-08-21 20:56:00"
$zhHandbackDt = "2016-08-21 20:56:27"
$deHandoffDt  = "2016-08-21 20:56:08"
$deHandbackDt = "2016-08-21 20:56:33"
$statusInSync = "Handed back: in sync with en-US"

$srcCommit   = "c1a9e6c2d4f8b3a07e5d1f29c8b6a4e3d7f90123"
$zhCnCommit  = "4f1a8c2e9b7d3f6a01c5e8b2d9f4a7c3e6b10945"
$deDeCommit  = "2e9b4f7a1c8d3e6f05a2c9b7d4f1a8e3c6b90247"

# ------------------------------------------------------------------
# Overview sheet (sheet1)
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "'$mdName"
$wsOverview.Range("C4").Value = "'.md"
$wsOverview.Range("E4").Value = "'$statusInSync"
$wsOverview.Range("F4").Value = "'$statusInSync"
$wsOverview.Range("G4").Value = "'$deHandoffDt"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/$mdPath",
    "",
    "",
    $mdPath
) | Out-Null

# ------------------------------------------------------------------
# zh-cn sheet (sheet2)
# ------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = "'$mdName"
$wsZh.Range("B4").Value = "'.md"
$wsZh.Range("C4").Value = "'$statusInSync"
$wsZh.Range("D4").Value = "'e2e"
$wsZh.Range("E4").Value = "'ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = "'$zhXlf"
$wsZh.Range("H4").Value = "'$zhHandoffDt"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = "'$mdName"
$wsZh.Range("J4").Value = "'$zhXlf"
$wsZh.Range("K4").Value = "'$zhHandbackDt"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/$mdPath",
    "",
    "",
    $mdName
) | Out-Null

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$zhCnCommit/$mdPath",
    "",
    "",
    $mdName
) | Out-Null

# ------------------------------------------------------------------
# de-de sheet (sheet3)
# ------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = "'$mdName"
$wsDe.Range("B4").Value = "'.md"
$wsDe.Range("C4").Value = "'$statusInSync"
$wsDe.Range("D4").Value = "'e2e"
$wsDe.Range("E4").Value = "'ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = "'$deXlf"
$wsDe.Range("H4").Value = "'$deHandoffDt"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = "'$mdName"
$wsDe.Range("J4").Value = "'$deXlf"
$wsDe.Range("K4").Value = "'$deHandbackDt"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/$mdPath",
    "",
    "",
    $mdName
) | Out-Null

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$deDeCommit/$mdPath",
    "",
    "",
    $mdName
) | Out-Null
